$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 3 with "w" in columns A through Q (matches the existing
# value used across row 2), extending the used range to A1:Q3.
$ws.Range("A3:Q3").Value = "w"
